# Generate Report for Handoff
#
# - Status text "Handed back: in sync with en-US" -> "Ready for handoff"
#   on the Overview sheet (E2, F2) and on each language sheet's Status
#   column (C2).
# - The "Latest HO Xliff Generate Date" / per-language handoff timestamps
#   are bumped forward a bit to reflect the new handoff run.
# - The Status column got visually narrower now that the text is shorter,
#   so the (autofit) column width shrinks to match.

$wb = $excel.ActiveWorkbook

$overview = $wb.Sheets.Item("Overview")
$zhcn     = $wb.Sheets.Item("zh-cn")
$dede     = $wb.Sheets.Item("de-de")

$newStatus = "Ready for handoff"

# --- Overview sheet -------------------------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-08-20 07:03:25"

# --- zh-cn sheet ------------------------------------------------------
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-08-20 07:03:21"

# --- de-de sheet ------------------------------------------------------
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-08-20 07:03:25"

# --- Column widths: the Status columns re-autofit narrower now that
#     "Ready for handoff" is shorter than "Handed back: in sync with
#     en-US". ColumnWidth rounds to the nearest pixel column width, so
#     feed it a value that lands on the same pixel cell as the target.
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth = 16.33
$dede.Columns.Item(3).ColumnWidth = 16.33
